$wb = $excel.ActiveWorkbook

$data = @{
    2 = @{ "B"=13.27156410718331; "C"=9.625922670097504; "D"=3.575400863928837; "E"=16.58294073349141; "F"=20.83286531184558; "H"=7.344005520526261; "N"=15.75157488290232; "O"=18.29973775294018 }
    3 = @{ "B"=12.60019494913002; "C"=9.134130970340939; "D"=3.549436474145852; "E"=15.63564054750957; "F"=20.70145259169839; "H"=7.344005520526261; "N"=15.80189921768247; "O"=18.27781497472204 }
    4 = @{ "B"=12.17052012960297; "C"=8.816595894821873; "D"=3.533264112081607; "E"=15.02879000762039; "F"=20.62864604092896; "H"=7.344005520526261; "N"=15.83465806696178; "O"=18.27090897552302 }
    5 = @{ "B"=11.99124259529006; "C"=8.683362882070044; "D"=3.526619002687049; "E"=14.77542321465228; "F"=20.60098457902129; "H"=7.344005520526261; "N"=15.84847605398079; "O"=18.26974068951908 }
    6 = @{ "B"=11.96122788212026; "C"=8.661010710042822; "D"=3.525512372818075; "E"=14.7329940929246; "F"=20.59651335638493; "H"=7.344005520526261; "N"=15.85079884278382; "O"=18.26964600694011 }
    7 = @{ "B"=12.16811896194693; "C"=8.814814476277427; "D"=3.533174711337349; "E"=15.02539720001694; "F"=20.62826482820009; "H"=7.344005520526261; "N"=15.83484252308212; "O"=18.27088655959344 }
    8 = @{ "B"=13.04381155157603; "C"=9.459628927021098; "D"=3.566497970554614; "E"=16.2616876183642; "F"=20.78593527466846; "H"=7.344005520526261; "N"=15.76854131182965; "O"=18.29081754348761 }
    9 = @{ "B"=14.61513578123397; "C"=10.597723095996; "D"=3.629883714108387; "E"=18.59154783222799; "F"=21.15632047419537; "H"=7.344005520526261; "N"=15.6532389882968; "O"=18.38192759954019 }
    10 = @{ "B"=15.6726344556512; "C"=11.35391702605415; "D"=3.675080027525733; "E"=20.24807268840485; "F"=21.46368247961469; "H"=7.344005520526261; "N"=15.57744171402299; "O"=18.48048428928114 }
    11 = @{ "B"=16.13143972783831; "C"=11.68014848791492; "D"=3.695305875644059; "E"=20.95931675908041; "F"=21.61064990407018; "H"=7.344005520526261; "N"=15.5448843084209; "O"=18.53212519522911 }
    12 = @{ "B"=16.30189895022955; "C"=11.80110409868711; "D"=3.702913747811815; "E"=21.22259862994759; "F"=21.66728315563908; "H"=7.344005520526261; "N"=15.53283142729869; "O"=18.55265092800852 }
    13 = @{ "B"=16.26533451519838; "C"=11.77516929915695; "D"=3.701277584070588; "E"=21.16616473536934; "F"=21.65504337810581; "H"=7.344005520526261; "N"=15.53541497004315; "O"=18.54818733690996 }
    14 = @{ "B"=16.14552965839578; "C"=11.69015143932472; "D"=3.695932824254285; "E"=20.98109818846065; "F"=21.61528980001895; "H"=7.344005520526261; "N"=15.54388718474472; "O"=18.53379446519936 }
    15 = @{ "B"=16.0717164097328; "C"=11.63773871956856; "D"=3.692652241612172; "E"=20.86695260535277; "F"=21.59106573389948; "H"=7.344005520526261; "N"=15.54911257279558; "O"=18.52510451503666 }
    16 = @{ "B"=15.6421954096912; "C"=11.33223713789318; "D"=3.673751224557625; "E"=20.2007426243093; "F"=21.45421759340264; "H"=7.344005520526261; "N"=15.57960805995792; "O"=18.47724573371858 }
    17 = @{ "B"=15.37293592415137; "C"=11.14025079965935; "D"=3.662068202999523; "E"=19.78123040252719; "F"=21.37206294210545; "H"=7.344005520526261; "N"=15.59880815061481; "O"=18.44962397603635 }
    18 = @{ "B"=15.21597382962659; "C"=11.02815435598203; "D"=3.655317169047275; "E"=19.53595291292994; "F"=21.32548549387339; "H"=7.344005520526261; "N"=15.61003258807859; "O"=18.43437799929937 }
    19 = @{ "B"=15.16247269633603; "C"=10.98991434969025; "D"=3.653026111751265; "E"=19.45222066482444; "F"=21.30983264517468; "H"=7.344005520526261; "N"=15.61386410593209; "O"=18.42932634707763 }
    20 = @{ "B"=15.40181611141165; "C"=11.16086122572586; "D"=3.663315139925944; "E"=19.82630039271546; "F"=21.38073884805473; "H"=7.344005520526261; "N"=15.59674553471742; "O"=18.45249803935712 }
    21 = @{ "B"=16.18080883948502; "C"=11.71519347866379; "D"=3.697504124719239; "E"=21.03562066414173; "F"=21.62694018897354; "H"=7.344005520526261; "N"=15.54139120701363; "O"=18.53799573908935 }
    22 = @{ "B"=16.67077820928148; "C"=12.06242850951304; "D"=3.719548642215879; "E"=21.79072833981449; "F"=21.79353469548956; "H"=7.344005520526261; "N"=15.50682184279566; "O"=18.59952471419356 }
    23 = @{ "B"=16.41104692461807; "C"=11.87848722027909; "D"=3.70781154451398; "E"=21.39092751753035; "F"=21.70411604581724; "H"=7.344005520526261; "N"=15.52512524655821; "O"=18.56617170026103 }
    24 = @{ "B"=15.3887661042717; "C"=11.15154861042544; "D"=3.6627515063501; "E"=19.80593700602815; "F"=21.37681442910545; "H"=7.344005520526261; "N"=15.59767746377773; "O"=18.45119669941386 }
    25 = @{ "B"=14.20659368664019; "C"=10.30371099249511; "D"=3.612964612971913; "E"=17.9437629592119; "F"=21.04976934090907; "H"=7.344005520526261; "N"=15.68286184824613; "O"=18.35171279602193 }
}

$ws = $wb.ActiveSheet

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
